$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 from "male" text to numeric 999
$ws.Range("A2").Value = 999

# Update B2 value
$ws.Range("B2").Value = -0.369

# Row 3 (female) stays the same gender label, update woe value
$ws.Range("A3").Value = "female"
$ws.Range("B3").Value = 0.266

# Add new row 4 with "male" and woe value 0.272
$ws.Range("A4").Value = "male"
$ws.Range("B4").Value = 0.272
